$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.043.03'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '3.382.71'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.48'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.51'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '3.380.24'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  +1.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('E11').Value = '  -3.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.382'
$ws.Range('E12').Value = '  -2.66%  '
$ws.Range('D13').Value = '3.957.31'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('E14').Value = '  +0.76%  '
$ws.Range('E15').Value = '  -3.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.66'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '3.379.70'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '61.185.87'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.84'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.76'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '376.97'
$ws.Range('E22').Value = '  -3.46%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.552'
$ws.Range('E23').Value = '  -2.44%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.519.67'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000127'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.96'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.181'
$ws.Range('E28').Value = '  +12.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.66'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.42'
$ws.Range('E31').Value = '  -3.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.08'
$ws.Range('E32').Value = '  -2.06%  '
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.44'
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('E36').Value = '  -4.71%  '
$ws.Range('E37').Value = '  -2.90%  '
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.71'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0757'
$ws.Range('E40').Value = '  -4.03%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.63'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.774'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('D47').Value = '2.546.05'
$ws.Range('E47').Value = '  +8.41%  '
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.96'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('E50').Value = '  +4.45%  '
$ws.Range('E51').Value = '  -1.89%  '
